$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.991.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.008.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.76%  "
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0777"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.304.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.734"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.025.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.913.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0809"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "221.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("E25").Value = "  -6.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("E27").Value = "  -5.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("E30").Value = "  -6.40%  "
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.453.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("E41").Value = "  -5.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.35%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.194.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.91%  "
